$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-10-22 06:36:20"

for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
